$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "2025/12/04 07:00"
$ws.Range("B39").Value = "32,744位本"
$ws.Range("C39").Value = "87位 広告・宣伝 (本)"
$ws.Range("D39").Value = "140位商業デザイン"
$ws.Range("E39").Value = "1,749位ビジネス実用本"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
